$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: 1/8 final match #1 ---
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2021-07-06"
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").Value = "Italy"
$ws.Range("D2").Value = "Spain"
$ws.Range("E2").Value = 0.3961661341853035
$ws.Range("F2").Value = 0.3067092651757188
$ws.Range("G2").Value = 0.2971246006389776
$ws.Range("H2").Value = 0.3521476748313809
$ws.Range("I2").Value = 0.3169329073482428

# --- Row 3: 1/8 final match #2 ---
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2021-07-07"
$ws.Range("B3").Style = "Normal"

$ws.Range("C3").Value = "England"
$ws.Range("D3").Value = "Denmark"
$ws.Range("E3").Value = 0.5524619634322795
$ws.Range("F3").Value = 0.1820493708072083
$ws.Range("G3").Value = 0.2654886657605121
$ws.Range("H3").Value = 0.4778795983689218
$ws.Range("I3").Value = 0.2207296066369154

# Remove the rows that are no longer part of the dataset (rows 4-9)
$ws.Range("A4:I9").Delete()
